$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 38
$ws.Cells.Item(38, 8).Value = 2637.5715
$ws.Cells.Item(38, 9).Value = 2495.2
$ws.Cells.Item(38, 10).Value = 2993.5
$ws.Cells.Item(38, 11).Value = 7485.599999999999
$ws.Cells.Item(38, 12).Value = 8980.5
$ws.Cells.Item(38, 13).Value = -7113.599999999999
$ws.Cells.Item(38, 14).Value = -9724.5
# Row 58
$ws.Cells.Item(58, 8).Value = 2928.111
$ws.Cells.Item(58, 9).Value = 58.833332
$ws.Cells.Item(58, 10).Value = 8666.666999999999
$ws.Cells.Item(58, 11).Value = 176.499996
$ws.Cells.Item(58, 12).Value = 26000.001
$ws.Cells.Item(58, 13).Value = -26.49999600000001
$ws.Cells.Item(58, 14).Value = -26300.001
# Row 64
$ws.Cells.Item(64, 8).Value = 7427.2856
$ws.Cells.Item(64, 9).Value = 7123
$ws.Cells.Item(64, 11).Value = 7123
$ws.Cells.Item(64, 13).Value = -6875
# Row 67
$ws.Cells.Item(67, 8).Value = 7427.2856
$ws.Cells.Item(67, 9).Value = 7123
$ws.Cells.Item(67, 11).Value = 7123
$ws.Cells.Item(67, 13).Value = -6265
# Row 76
$ws.Cells.Item(76, 8).Value = 6004
$ws.Cells.Item(76, 9).Value = 0
$ws.Cells.Item(76, 10).Value = 6004
$ws.Cells.Item(76, 11).Value = 0
$ws.Cells.Item(76, 12).Value = 6004
$ws.Cells.Item(76, 13).ClearContents()
$ws.Cells.Item(76, 14).Value = -6634
# Row 79
$ws.Cells.Item(79, 8).Value = 6004
$ws.Cells.Item(79, 9).Value = 0
$ws.Cells.Item(79, 10).Value = 6004
$ws.Cells.Item(79, 11).Value = 0
$ws.Cells.Item(79, 12).Value = 6004
$ws.Cells.Item(79, 13).ClearContents()
$ws.Cells.Item(79, 14).Value = -8188
# Row 86
$ws.Cells.Item(86, 8).Value = 8259.684999999999
$ws.Cells.Item(86, 9).Value = 9123.625
$ws.Cells.Item(86, 10).Value = 7631.364
$ws.Cells.Item(86, 11).Value = 9123.625
$ws.Cells.Item(86, 12).Value = 7631.364
$ws.Cells.Item(86, 13).Value = -8000.625
$ws.Cells.Item(86, 14).Value = -9877.364
# Row 89
$ws.Cells.Item(89, 8).Value = 8259.684999999999
$ws.Cells.Item(89, 9).Value = 9123.625
$ws.Cells.Item(89, 10).Value = 7631.364
$ws.Cells.Item(89, 11).Value = 45618.125
$ws.Cells.Item(89, 12).Value = 38156.82
$ws.Cells.Item(89, 13).Value = -40002.125
$ws.Cells.Item(89, 14).Value = -49388.82
# Row 138
$ws.Cells.Item(138, 8).Value = 3859.8462
$ws.Cells.Item(138, 9).Value = 2669.3076
$ws.Cells.Item(138, 10).Value = 4455.115
$ws.Cells.Item(138, 11).Value = 8007.9228
$ws.Cells.Item(138, 12).Value = 13365.345
$ws.Cells.Item(138, 13).Value = -2867.9228
$ws.Cells.Item(138, 14).Value = -23645.345

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 4948.97
$ws.Cells.Item(32, 9).Value = 1952.931
$ws.Cells.Item(32, 10).Value = 24999.385
$ws.Cells.Item(32, 11).Value = 1952.931
$ws.Cells.Item(32, 12).Value = 24999.385
$ws.Cells.Item(32, 13).Value = -1665.931
$ws.Cells.Item(32, 14).Value = -25573.385
# Row 61
$ws.Cells.Item(61, 8).Value = 4351732
$ws.Cells.Item(61, 9).Value = 3941.2896
$ws.Cells.Item(61, 11).Value = 3941.2896
$ws.Cells.Item(61, 13).Value = -3729.2896
# Row 132
$ws.Cells.Item(132, 8).Value = 965738.9
$ws.Cells.Item(132, 9).Value = 1091144
$ws.Cells.Item(132, 11).Value = 3273432
$ws.Cells.Item(132, 13).Value = -3270902
# Row 136
$ws.Cells.Item(136, 8).Value = 4351732
$ws.Cells.Item(136, 9).Value = 3941.2896
$ws.Cells.Item(136, 11).Value = 11823.8688
$ws.Cells.Item(136, 13).Value = -9273.8688

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Cells.Item(134, 8).Value = 4394375
$ws.Cells.Item(134, 9).Value = 6957.9443
$ws.Cells.Item(134, 11).Value = 20873.8329
$ws.Cells.Item(134, 13).Value = -18338.8329

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 2060873.2
$ws.Cells.Item(31, 9).Value = 3269854.8
$ws.Cells.Item(31, 10).Value = 5604.9
$ws.Cells.Item(31, 11).Value = 3269854.8
$ws.Cells.Item(31, 12).Value = 5604.9
$ws.Cells.Item(31, 13).Value = -3269559.8
$ws.Cells.Item(31, 14).Value = -6194.9
# Row 34
$ws.Cells.Item(34, 8).Value = 2060873.2
$ws.Cells.Item(34, 9).Value = 3269854.8
$ws.Cells.Item(34, 10).Value = 5604.9
$ws.Cells.Item(34, 11).Value = 3269854.8
$ws.Cells.Item(34, 12).Value = 5604.9
$ws.Cells.Item(34, 13).Value = -3269652.8
$ws.Cells.Item(34, 14).Value = -6008.9
# Row 105
$ws.Cells.Item(105, 8).Value = 9554.532999999999
$ws.Cells.Item(105, 10).Value = 5710.4287
$ws.Cells.Item(105, 12).Value = 5710.4287
$ws.Cells.Item(105, 14).Value = -9204.4287
# Row 107
$ws.Cells.Item(107, 8).Value = 1479.8
$ws.Cells.Item(107, 9).Value = 533
$ws.Cells.Item(107, 11).Value = 533
$ws.Cells.Item(107, 13).Value = 1387
# Row 108
$ws.Cells.Item(108, 8).Value = 46809.855
$ws.Cells.Item(108, 10).Value = 46809.855
$ws.Cells.Item(108, 12).Value = 46809.855
$ws.Cells.Item(108, 14).Value = -54489.855
# Row 109
$ws.Cells.Item(109, 8).Value = 59999
$ws.Cells.Item(109, 10).Value = 59999
$ws.Cells.Item(109, 12).Value = 59999
$ws.Cells.Item(109, 14).Value = -62079
# Row 132
$ws.Cells.Item(132, 8).Value = 8034.9165
$ws.Cells.Item(132, 9).Value = 3875.7368
$ws.Cells.Item(132, 11).Value = 11627.2104
$ws.Cells.Item(132, 13).Value = -9097.2104
# Row 141
$ws.Cells.Item(141, 8).Value = 292613.06
$ws.Cells.Item(141, 10).Value = 327633.72
$ws.Cells.Item(141, 12).Value = 327633.72
$ws.Cells.Item(141, 14).Value = -337993.72

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Cells.Item(2, 8).Value = 152.90909
$ws.Cells.Item(2, 9).Value = 87.46154
$ws.Cells.Item(2, 10).Value = 247.44444
$ws.Cells.Item(2, 11).Value = 524.76924
$ws.Cells.Item(2, 12).Value = 1484.66664
$ws.Cells.Item(2, 13).Value = -411.76924
$ws.Cells.Item(2, 14).Value = -1710.66664
# Row 80
$ws.Cells.Item(80, 8).Value = 15000
$ws.Cells.Item(80, 9).Value = 0
$ws.Cells.Item(80, 10).Value = 15000
$ws.Cells.Item(80, 11).Value = 0
$ws.Cells.Item(80, 12).Value = 45000
$ws.Cells.Item(80, 13).ClearContents()
$ws.Cells.Item(80, 14).Value = -46872
# Row 83
$ws.Cells.Item(83, 8).Value = 15000
$ws.Cells.Item(83, 9).Value = 0
$ws.Cells.Item(83, 10).Value = 15000
$ws.Cells.Item(83, 11).Value = 0
$ws.Cells.Item(83, 12).Value = 135000
$ws.Cells.Item(83, 13).ClearContents()
$ws.Cells.Item(83, 14).Value = -144360
# Row 113
$ws.Cells.Item(113, 8).Value = 1859.8334
$ws.Cells.Item(113, 9).Value = 1098.6666
$ws.Cells.Item(113, 11).Value = 3295.9998
$ws.Cells.Item(113, 13).Value = -1125.9998
# Row 131
$ws.Cells.Item(131, 8).Value = 3992.257
$ws.Cells.Item(131, 10).Value = 5146.5
$ws.Cells.Item(131, 12).Value = 15439.5
$ws.Cells.Item(131, 14).Value = -25519.5
# Row 132
$ws.Cells.Item(132, 8).Value = 5664.6924
$ws.Cells.Item(132, 9).Value = 1566.6666
$ws.Cells.Item(132, 10).Value = 6199.2173
$ws.Cells.Item(132, 11).Value = 14099.9994
$ws.Cells.Item(132, 12).Value = 55792.95570000001
$ws.Cells.Item(132, 13).Value = -11569.9994
$ws.Cells.Item(132, 14).Value = -60852.95570000001

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Cells.Item(126, 8).Value = 13564.1875
$ws.Cells.Item(126, 9).Value = 17911.545
$ws.Cells.Item(126, 10).Value = 4000
$ws.Cells.Item(126, 11).Value = 53734.63499999999
$ws.Cells.Item(126, 12).Value = 12000
$ws.Cells.Item(126, 13).Value = -51264.63499999999
$ws.Cells.Item(126, 14).Value = -16940
# Row 132
$ws.Cells.Item(132, 8).Value = 21775.133
$ws.Cells.Item(132, 9).Value = 24341.125
$ws.Cells.Item(132, 10).Value = 18842.572
$ws.Cells.Item(132, 11).Value = 73023.375
$ws.Cells.Item(132, 12).Value = 56527.716
$ws.Cells.Item(132, 13).Value = -70493.375
$ws.Cells.Item(132, 14).Value = -61587.716

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Cells.Item(61, 8).Value = 13369.5
$ws.Cells.Item(61, 9).Value = 14129.846
$ws.Cells.Item(61, 10).Value = 3485
$ws.Cells.Item(61, 11).Value = 14129.846
$ws.Cells.Item(61, 12).Value = 3485
$ws.Cells.Item(61, 13).Value = -13927.846
$ws.Cells.Item(61, 14).Value = -3889
# Row 113
$ws.Cells.Item(113, 8).Value = 13369.5
$ws.Cells.Item(113, 9).Value = 14129.846
$ws.Cells.Item(113, 10).Value = 3485
$ws.Cells.Item(113, 11).Value = 14129.846
$ws.Cells.Item(113, 12).Value = 3485
$ws.Cells.Item(113, 13).Value = -11959.846
$ws.Cells.Item(113, 14).Value = -7825
# Row 132
$ws.Cells.Item(132, 8).Value = 3772806
$ws.Cells.Item(132, 9).Value = 9740643
$ws.Cells.Item(132, 10).Value = 3645.9473
$ws.Cells.Item(132, 11).Value = 29221929
$ws.Cells.Item(132, 12).Value = 10937.8419
$ws.Cells.Item(132, 13).Value = -29219399
$ws.Cells.Item(132, 14).Value = -15997.8419

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Cells.Item(132, 8).Value = 7249170
$ws.Cells.Item(132, 9).Value = 9806730
$ws.Cells.Item(132, 10).Value = 2750
$ws.Cells.Item(132, 11).Value = 29420190
$ws.Cells.Item(132, 12).Value = 8250
$ws.Cells.Item(132, 13).Value = -29417660
$ws.Cells.Item(132, 14).Value = -13310
